{"js": "// Replace the customer name, address, invoice number, and invoice date\n// throughout the document (Safety Certificate body text, the duplicate\n// summary paragraph, and the serial-number certification paragraph).\nconst replacements = [\n  { from: \"Prem Singh\", to: \"ABC\" },\n  { from: \"dsfghbujikoml\", to: \"Strret no 1, Kaithal\" },\n  { from: \"ME/2025-26/359\", to: \"ME/2025-26/25\" },\n  { from: \"05-04-2025\", to: \"20-04-2025\" },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the customer name, address, invoice number, and invoice date\n# throughout the document (Safety Certificate body text, the duplicate\n# summary paragraph, and the serial-number certification paragraph).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ From = \"Prem Singh\"; To = \"ABC\" },\n    @{ From = \"dsfghbujikoml\"; To = \"Strret no 1, Kaithal\" },\n    @{ From = \"ME/2025-26/359\"; To = \"ME/2025-26/25\" },\n    @{ From = \"05-04-2025\"; To = \"20-04-2025\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.From\n    $find.Replacement.Text = $r.To\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$find.Text, [ref]$false, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null\n}\n"}
